$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = "NSE:APLAPOLLO"
$ws.Range("C2").Value2 = "NSE:AARVI"
$ws.Range("D2").Value2 = "NSE:ABB"
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value2 = "NSE:LICHSGFIN"

$ws.Range("B3").Value2 = "NSE:ASTRAL"
$ws.Range("C3").Value2 = "NSE:ALLCARGO"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

$ws.Range("B4").Value2 = "NSE:AUROPHARMA"
$ws.Range("C4").Value2 = "NSE:ALLSEC"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

$ws.Range("B5").Value2 = "NSE:AVANTIFEED"
$ws.Range("C5").Value2 = "NSE:ASPINWALL"
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()

$ws.Range("B6").Value2 = "NSE:AWL"
$ws.Range("C6").Value2 = "NSE:ASTERDM"
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").ClearContents()

$ws.Range("B7").Value2 = "NSE:CMSINFO"
$ws.Range("C7").Value2 = "NSE:AUTOIND"
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()

$ws.Range("B8").Value2 = "NSE:FEDERALBNK"
$ws.Range("C8").Value2 = "NSE:BANG"
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

$ws.Range("B9").Value2 = "NSE:GODREJIND"
$ws.Range("C9").Value2 = "NSE:BEARDSELL"
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()

$ws.Range("B10").Value2 = "NSE:GROBTEA"
$ws.Range("C10").Value2 = "NSE:CHALET"
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("B11").Value2 = "NSE:HGS"
$ws.Range("C11").Value2 = "NSE:CHENNPETRO"
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

$ws.Range("B12").Value2 = "NSE:KOTHARIPRO"
$ws.Range("C12").Value2 = "NSE:CONFIPET"
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()

$ws.Range("B13").Value2 = "NSE:LICHSGFIN"
$ws.Range("C13").Value2 = "NSE:DEVYANI"
$ws.Range("D13").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()

$ws.Range("B14").Value2 = "NSE:MAKEINDIA"
$ws.Range("C14").Value2 = "NSE:DMCC"
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()

$ws.Range("B15").Value2 = "NSE:MANYAVAR"
$ws.Range("C15").Value2 = "NSE:GENSOL"
$ws.Range("D15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()

$ws.Range("B16").Value2 = "NSE:MIRCELECTR"
$ws.Range("C16").Value2 = "NSE:GMRP&UI"
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()

$ws.Range("B17").Value2 = "NSE:MOQUALITY"
$ws.Range("C17").Value2 = "NSE:INOXWIND"
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("F17").ClearContents()

$ws.Range("B18").Value2 = "NSE:NEWGEN"
$ws.Range("C18").Value2 = "NSE:IOLCP"
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").ClearContents()

$ws.Range("B19").Value2 = "NSE:PGHL"
$ws.Range("C19").Value2 = "NSE:KABRAEXTRU"
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("B20").Value2 = "NSE:PNB"
$ws.Range("C20").Value2 = "NSE:KMSUGAR"
$ws.Range("D20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()

$ws.Range("B21").Value2 = "NSE:RAMRAT"
$ws.Range("C21").Value2 = "NSE:KSB"
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("F21").ClearContents()

$ws.Range("B22").ClearContents()
$ws.Range("C22").Value2 = "NSE:LIBERTSHOE"
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("F22").ClearContents()

$ws.Range("B23").ClearContents()
$ws.Range("C23").Value2 = "NSE:LTTS"
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").ClearContents()

$ws.Range("B24").ClearContents()
$ws.Range("C24").Value2 = "NSE:MALLCOM"
$ws.Range("D24").ClearContents()
$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()

$ws.Range("B25").ClearContents()
$ws.Range("C25").Value2 = "NSE:MANALIPETC"
$ws.Range("D25").ClearContents()
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()

$ws.Range("B26").ClearContents()
$ws.Range("C26").Value2 = "NSE:MAZDOCK"
$ws.Range("D26").ClearContents()
$ws.Range("E26").ClearContents()
$ws.Range("F26").ClearContents()

$ws.Range("B27").ClearContents()
$ws.Range("C27").Value2 = "NSE:MBAPL"
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("F27").ClearContents()

$ws.Range("B28").ClearContents()
$ws.Range("C28").Value2 = "NSE:MOIL"
$ws.Range("D28").ClearContents()
$ws.Range("E28").ClearContents()
$ws.Range("F28").ClearContents()

$ws.Range("B29").ClearContents()
$ws.Range("C29").Value2 = "NSE:OLECTRA"
$ws.Range("D29").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("F29").ClearContents()

$ws.Range("B30").ClearContents()
$ws.Range("C30").Value2 = "NSE:ORIENTBELL"
$ws.Range("D30").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("F30").ClearContents()

$ws.Range("B31").ClearContents()
$ws.Range("C31").Value2 = "NSE:PAISALO"
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("F31").ClearContents()

$ws.Range("B32").ClearContents()
$ws.Range("C32").Value2 = "NSE:PANACEABIO"
$ws.Range("D32").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("F32").ClearContents()

$ws.Range("B33").ClearContents()
$ws.Range("C33").Value2 = "NSE:PARAS"
$ws.Range("D33").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("F33").ClearContents()

$ws.Range("B34").ClearContents()
$ws.Range("C34").Value2 = "NSE:PIXTRANS"
$ws.Range("D34").ClearContents()
$ws.Range("E34").ClearContents()
$ws.Range("F34").ClearContents()

$ws.Range("B35").ClearContents()
$ws.Range("C35").Value2 = "NSE:POLYPLEX"
$ws.Range("D35").ClearContents()
$ws.Range("E35").ClearContents()
$ws.Range("F35").ClearContents()

$ws.Range("B36").ClearContents()
$ws.Range("C36").Value2 = "NSE:PYRAMID"
$ws.Range("D36").ClearContents()
$ws.Range("E36").ClearContents()
$ws.Range("F36").ClearContents()

$ws.Range("B37").ClearContents()
$ws.Range("C37").Value2 = "NSE:QUICKHEAL"
$ws.Range("D37").ClearContents()
$ws.Range("E37").ClearContents()
$ws.Range("F37").ClearContents()

$ws.Range("B38").ClearContents()
$ws.Range("C38").Value2 = "NSE:RAMAPHO"
$ws.Range("D38").ClearContents()
$ws.Range("E38").ClearContents()
$ws.Range("F38").ClearContents()

$ws.Rows.Item(39).Delete()

